$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": refresh scraped_at timestamps in column K (rows 2-52) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 11).Value = "2025-11-07T02:48:24.023086+00:00"
$ws1.Cells.Item(3, 11).Value = "2025-11-07T02:48:24.023108+00:00"
$ws1.Cells.Item(4, 11).Value = "2025-11-07T02:48:24.023118+00:00"
$ws1.Cells.Item(5, 11).Value = "2025-11-07T02:48:26.025171+00:00"
$ws1.Cells.Item(6, 11).Value = "2025-11-07T02:48:26.025250+00:00"
$ws1.Cells.Item(7, 11).Value = "2025-11-07T02:48:26.025287+00:00"
$ws1.Cells.Item(8, 11).Value = "2025-11-07T02:48:28.469890+00:00"
$ws1.Cells.Item(9, 11).Value = "2025-11-07T02:48:30.833939+00:00"
$ws1.Cells.Item(10, 11).Value = "2025-11-07T02:48:30.833956+00:00"
$ws1.Cells.Item(11, 11).Value = "2025-11-07T02:48:30.833965+00:00"
$ws1.Cells.Item(12, 11).Value = "2025-11-07T02:48:33.243613+00:00"
$ws1.Cells.Item(13, 11).Value = "2025-11-07T02:48:33.243637+00:00"
$ws1.Cells.Item(14, 11).Value = "2025-11-07T02:48:33.243647+00:00"
$ws1.Cells.Item(15, 11).Value = "2025-11-07T02:48:33.243656+00:00"
$ws1.Cells.Item(16, 11).Value = "2025-11-07T02:48:38.123709+00:00"
$ws1.Cells.Item(17, 11).Value = "2025-11-07T02:48:40.615052+00:00"
$ws1.Cells.Item(18, 11).Value = "2025-11-07T02:48:43.101589+00:00"
$ws1.Cells.Item(19, 11).Value = "2025-11-07T02:48:43.101609+00:00"
$ws1.Cells.Item(20, 11).Value = "2025-11-07T02:48:45.558452+00:00"
$ws1.Cells.Item(21, 11).Value = "2025-11-07T02:48:47.964830+00:00"
$ws1.Cells.Item(22, 11).Value = "2025-11-07T02:48:47.964850+00:00"
$ws1.Cells.Item(23, 11).Value = "2025-11-07T02:48:47.964858+00:00"
$ws1.Cells.Item(24, 11).Value = "2025-11-07T02:48:50.088811+00:00"
$ws1.Cells.Item(25, 11).Value = "2025-11-07T02:48:50.088843+00:00"
$ws1.Cells.Item(26, 11).Value = "2025-11-07T02:48:50.088864+00:00"
$ws1.Cells.Item(27, 11).Value = "2025-11-07T02:48:52.541358+00:00"
$ws1.Cells.Item(28, 11).Value = "2025-11-07T02:48:52.541391+00:00"
$ws1.Cells.Item(29, 11).Value = "2025-11-07T02:48:52.541412+00:00"
$ws1.Cells.Item(30, 11).Value = "2025-11-07T02:48:52.541431+00:00"
$ws1.Cells.Item(31, 11).Value = "2025-11-07T02:48:52.541448+00:00"
$ws1.Cells.Item(32, 11).Value = "2025-11-07T02:48:55.059239+00:00"
$ws1.Cells.Item(33, 11).Value = "2025-11-07T02:48:57.499048+00:00"
$ws1.Cells.Item(34, 11).Value = "2025-11-07T02:48:57.499083+00:00"
$ws1.Cells.Item(35, 11).Value = "2025-11-07T02:48:57.499106+00:00"
$ws1.Cells.Item(36, 11).Value = "2025-11-07T02:48:59.959114+00:00"
$ws1.Cells.Item(37, 11).Value = "2025-11-07T02:48:59.959146+00:00"
$ws1.Cells.Item(38, 11).Value = "2025-11-07T02:48:59.959166+00:00"
$ws1.Cells.Item(39, 11).Value = "2025-11-07T02:49:01.984797+00:00"
$ws1.Cells.Item(40, 11).Value = "2025-11-07T02:49:01.984832+00:00"
$ws1.Cells.Item(41, 11).Value = "2025-11-07T02:49:01.984852+00:00"
$ws1.Cells.Item(42, 11).Value = "2025-11-07T02:49:01.984870+00:00"
$ws1.Cells.Item(43, 11).Value = "2025-11-07T02:49:01.984889+00:00"
$ws1.Cells.Item(44, 11).Value = "2025-11-07T02:49:01.984906+00:00"
$ws1.Cells.Item(45, 11).Value = "2025-11-07T02:49:01.984922+00:00"
$ws1.Cells.Item(46, 11).Value = "2025-11-07T02:49:01.984938+00:00"
$ws1.Cells.Item(47, 11).Value = "2025-11-07T02:49:04.097065+00:00"
$ws1.Cells.Item(48, 11).Value = "2025-11-07T02:49:04.097113+00:00"
$ws1.Cells.Item(49, 11).Value = "2025-11-07T02:49:08.175464+00:00"
$ws1.Cells.Item(50, 11).Value = "2025-11-07T02:49:08.175495+00:00"
$ws1.Cells.Item(51, 11).Value = "2025-11-07T02:49:10.212178+00:00"
$ws1.Cells.Item(52, 11).Value = "2025-11-07T02:49:10.212199+00:00"

# --- Sheet "returned": remove the single (now stale) data row ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()

# --- Sheet "new_injured": remove the single (now stale) data row ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Delete()
